$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set B7 (Experimental value) to the text string "true" without Excel's
# automatic boolean-literal coercion (which would store it as t="b" and
# also avoid triggering a quote-prefix style). We build it via a formula
# that evaluates to the literal text "true", then convert the formula
# result to a static value in-place (Copy + PasteSpecial values-only),
# which preserves the existing cell style (s="2").
$ws.Range("B7").Formula = '=LOWER(TEXT(TRUE,""))'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Update the Date value to the new timestamp.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
